$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values for rows 2-6 (columns E-T).
# Columns A-D (Sending cluster / Ligand symbol / Receptor symbol / Target cluster)
# keep the same text content; only the numeric expression/specificity values change.

$data = @{
    2 = @{ E=3; F=1; G=0.7878926666666667; H=2.363678;
           M=1.325636333333333; N=3.976909;
           O=0.1236745921078505; P=0.1236745921078504;
           Q=1.044459145700222;  R=9.400132311302;
           S=0.1236745921078505; T=0.1236745921078504 }
    3 = @{ E=3; F=1; G=0.7878926666666667; H=2.363678;
           O=0.1565567708132977; P=0.1565567708132977;
           Q=1.322156380792;     R=11.899407427128;
           S=0.1565567708132977; T=0.1565567708132977 }
    4 = @{ E=3; F=1; G=0.7878926666666667; H=2.363678;
           M=2.239683666666667; N=6.719051;
           O=0.2089501901544251; P=0.2089501901544251;
           Q=1.764630336619778;  R=15.881673029578;
           S=0.2089501901544251; T=0.2089501901544251 }
    5 = @{ E=3; F=1; G=0.7878926666666667; H=2.363678;
           M=2.783146333333333; N=8.349439;
           O=0.2596522733244283; P=0.2596522733244283;
           Q=2.192820586293556;  R=19.735385276642;
           S=0.2596522733244283; T=0.2596522733244283 }
    6 = @{ E=3; F=1; G=0.7878926666666667; H=2.363678;
           M=2.692186; N=8.076558;
           O=0.2511661735999985; P=0.2511661735999985;
           Q=2.121153606702667;  R=19.090382460324;
           S=0.2511661735999985; T=0.2511661735999985 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}
